$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, reuse the header style from E1 via copy/paste-format
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Per-row time_taken timestamps (F2:F77)
$timeTaken = @(
    "2021-10-05 13:42:19.040075",
    "2021-10-05 13:42:19.040088",
    "2021-10-05 13:42:19.040092",
    "2021-10-05 13:42:19.040096",
    "2021-10-05 13:42:19.040099",
    "2021-10-05 13:42:19.040102",
    "2021-10-05 13:42:19.040105",
    "2021-10-05 13:42:19.040109",
    "2021-10-05 13:42:19.040112",
    "2021-10-05 13:42:19.040115",
    "2021-10-05 13:42:19.040118",
    "2021-10-05 13:42:19.040121",
    "2021-10-05 13:42:19.040124",
    "2021-10-05 13:42:19.040127",
    "2021-10-05 13:42:19.040129",
    "2021-10-05 13:42:19.040132",
    "2021-10-05 13:42:19.040136",
    "2021-10-05 13:42:19.040139",
    "2021-10-05 13:42:19.040142",
    "2021-10-05 13:42:19.040145",
    "2021-10-05 13:42:19.040148",
    "2021-10-05 13:42:19.040151",
    "2021-10-05 13:42:19.040153",
    "2021-10-05 13:42:19.040156",
    "2021-10-05 13:42:19.040160",
    "2021-10-05 13:42:19.040163",
    "2021-10-05 13:42:19.040166",
    "2021-10-05 13:42:19.040168",
    "2021-10-05 13:42:19.040171",
    "2021-10-05 13:42:19.040174",
    "2021-10-05 13:42:19.040177",
    "2021-10-05 13:42:19.040180",
    "2021-10-05 13:42:19.040184",
    "2021-10-05 13:42:19.040187",
    "2021-10-05 13:42:19.040190",
    "2021-10-05 13:42:19.040193",
    "2021-10-05 13:42:19.040196",
    "2021-10-05 13:42:19.040198",
    "2021-10-05 13:42:19.040201",
    "2021-10-05 13:42:19.040204",
    "2021-10-05 13:42:19.040208",
    "2021-10-05 13:42:19.040211",
    "2021-10-05 13:42:19.040214",
    "2021-10-05 13:42:19.040217",
    "2021-10-05 13:42:19.040220",
    "2021-10-05 13:42:19.040223",
    "2021-10-05 13:42:19.040226",
    "2021-10-05 13:42:19.040229",
    "2021-10-05 13:42:19.040231",
    "2021-10-05 13:42:19.040234",
    "2021-10-05 13:42:19.040237",
    "2021-10-05 13:42:19.040240",
    "2021-10-05 13:42:19.040243",
    "2021-10-05 13:42:19.040247",
    "2021-10-05 13:42:19.040250",
    "2021-10-05 13:42:19.040253",
    "2021-10-05 13:42:19.040255",
    "2021-10-05 13:42:19.040258",
    "2021-10-05 13:42:19.040261",
    "2021-10-05 13:42:19.040264",
    "2021-10-05 13:42:19.040267",
    "2021-10-05 13:42:19.040270",
    "2021-10-05 13:42:19.040273",
    "2021-10-05 13:42:19.040276",
    "2021-10-05 13:42:19.040280",
    "2021-10-05 13:42:19.040284",
    "2021-10-05 13:42:19.040287",
    "2021-10-05 13:42:19.040290",
    "2021-10-05 13:42:19.040293",
    "2021-10-05 13:42:19.040296",
    "2021-10-05 13:42:19.040298",
    "2021-10-05 13:42:19.040301",
    "2021-10-05 13:42:19.040304",
    "2021-10-05 13:42:19.040307",
    "2021-10-05 13:42:19.040310",
    "2021-10-05 13:42:19.040313"
)

for ($i = 0; $i -lt $timeTaken.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

$excel.CutCopyMode = $false

